# Continuation de l'application sous classes
# Extend the participants table from row 24 (N_Ano=23) down to row 30 (N_Ano=29),
# mirroring the existing pattern: column A gets the sequential participant
# number (row-1), columns B:P are left blank placeholders for data entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastCol = 16   # column P

for ($row = 25; $row -le 30; $row++) {
    $ws.Cells.Item($row, 1).Value = $row - 1

    for ($col = 2; $col -le $lastCol; $col++) {
        $ws.Cells.Item($row, $col).Value = ""
    }
}
